$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 90.3
$ws.Range("I9").Value = 111.53846
$ws.Range("J9").Value = 50.857143
$ws.Range("K9").Value = 111.53846
$ws.Range("L9").Value = 50.857143
$ws.Range("M9").Value = 57.46154
$ws.Range("N9").Value = -388.857143
$ws.Range("H17").Value = 1991.24
$ws.Range("J17").Value = 1991.24
$ws.Range("L17").Value = 5973.72
$ws.Range("N17").Value = -6309.72
$ws.Range("H113").Value = 4566.1113
$ws.Range("I113").Value = 3549.1667
$ws.Range("J113").Value = 6600
$ws.Range("K113").Value = 3549.1667
$ws.Range("L113").Value = 6600
$ws.Range("M113").Value = -295.1667000000002
$ws.Range("N113").Value = -13108
$ws.Range("H125").Value = 6947551
$ws.Range("I125").Value = 1880.25
$ws.Range("J125").Value = 13893222
$ws.Range("K125").Value = 16922.25
$ws.Range("L125").Value = 125038998
$ws.Range("M125").Value = -14462.25
$ws.Range("N125").Value = -125043918
$ws.Range("H137").Value = 2260.0527
$ws.Range("J137").Value = 3821.4375
$ws.Range("L137").Value = 11464.3125
$ws.Range("N137").Value = -16564.3125
$ws.Range("H138").Value = 5116.564
$ws.Range("I138").Value = 2545.9583
$ws.Range("K138").Value = 7637.874899999999
$ws.Range("M138").Value = -2497.874899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2092.182
$ws.Range("I45").Value = 2101.4
$ws.Range("K45").Value = 2101.4
$ws.Range("M45").Value = -1724.4
$ws.Range("H61").Value = 4250.265
$ws.Range("I61").Value = 3734.9524
$ws.Range("J61").Value = 7342.143
$ws.Range("K61").Value = 3734.9524
$ws.Range("L61").Value = 7342.143
$ws.Range("M61").Value = -3522.9524
$ws.Range("N61").Value = -7766.143
$ws.Range("H136").Value = 4250.265
$ws.Range("I136").Value = 3734.9524
$ws.Range("J136").Value = 7342.143
$ws.Range("K136").Value = 11204.8572
$ws.Range("L136").Value = 22026.429
$ws.Range("M136").Value = -8654.8572
$ws.Range("N136").Value = -27126.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 26931.195
$ws.Range("I134").Value = 3995.1333
$ws.Range("J134").Value = 69936.31
$ws.Range("K134").Value = 11985.3999
$ws.Range("L134").Value = 209808.93
$ws.Range("M134").Value = -9450.3999
$ws.Range("N134").Value = -214878.93

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26064.238
$ws.Range("J31").Value = 37265.69
$ws.Range("L31").Value = 37265.69
$ws.Range("N31").Value = -37855.69
$ws.Range("H34").Value = 26064.238
$ws.Range("J34").Value = 37265.69
$ws.Range("L34").Value = 37265.69
$ws.Range("N34").Value = -37669.69
$ws.Range("H122").Value = 2666
$ws.Range("I122").Value = 1999.5
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 5998.5
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -3548.5
$ws.Range("N122").Value = -16897

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 22818.584
$ws.Range("I87").Value = 12364.6
$ws.Range("J87").Value = 30285.715
$ws.Range("K87").Value = 37093.8
$ws.Range("L87").Value = 90857.145
$ws.Range("M87").Value = -35845.8
$ws.Range("N87").Value = -93353.145
$ws.Range("H90").Value = 22818.584
$ws.Range("I90").Value = 12364.6
$ws.Range("J90").Value = 30285.715
$ws.Range("K90").Value = 111281.4
$ws.Range("L90").Value = 272571.435
$ws.Range("M90").Value = -105041.4
$ws.Range("N90").Value = -285051.435
$ws.Range("H107").Value = 74250.82000000001
$ws.Range("I107").Value = 1111.625
$ws.Range("J107").Value = 103506.5
$ws.Range("K107").Value = 3334.875
$ws.Range("L107").Value = 310519.5
$ws.Range("M107").Value = -1414.875
$ws.Range("N107").Value = -314359.5
$ws.Range("H113").Value = 1545402.1
$ws.Range("I113").Value = 6174092.5
$ws.Range("J113").Value = 2505.2222
$ws.Range("K113").Value = 18522277.5
$ws.Range("L113").Value = 7515.6666
$ws.Range("M113").Value = -18520107.5
$ws.Range("N113").Value = -11855.6666
$ws.Range("H122").Value = 127269.75
$ws.Range("J122").Value = 169278.17
$ws.Range("L122").Value = 1523503.53
$ws.Range("N122").Value = -1528403.53
$ws.Range("H130").Value = 27500.25
$ws.Range("I130").Value = 5001
$ws.Range("K130").Value = 15003
$ws.Range("M130").Value = -9983
$ws.Range("H131").Value = 2254.9177
$ws.Range("I131").Value = 1738.7273
$ws.Range("J131").Value = 2331.6487
$ws.Range("K131").Value = 5216.1819
$ws.Range("L131").Value = 6994.946100000001
$ws.Range("M131").Value = -176.1818999999996
$ws.Range("N131").Value = -17074.9461
$ws.Range("H132").Value = 582688.6
$ws.Range("I132").Value = 145695.42
$ws.Range("K132").Value = 1311258.78
$ws.Range("M132").Value = -1308728.78
$ws.Range("H139").Value = 5371.028
$ws.Range("I139").Value = 1692.1111
$ws.Range("K139").Value = 5076.3333
$ws.Range("M139").Value = 63.66669999999976

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 744614.25
$ws.Range("I80").Value = 774026
$ws.Range("J80").Value = 717303.4
$ws.Range("K80").Value = 774026
$ws.Range("L80").Value = 717303.4
$ws.Range("M80").Value = -773028
$ws.Range("N80").Value = -719299.4
$ws.Range("H83").Value = 744614.25
$ws.Range("I83").Value = 774026
$ws.Range("J83").Value = 717303.4
$ws.Range("K83").Value = 3870130
$ws.Range("L83").Value = 3586517
$ws.Range("M83").Value = -3865138
$ws.Range("N83").Value = -3596501
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 24943.354
$ws.Range("I132").Value = 4943.6387
$ws.Range("J132").Value = 84942.5
$ws.Range("K132").Value = 14830.9161
$ws.Range("L132").Value = 254827.5
$ws.Range("M132").Value = -12300.9161
$ws.Range("N132").Value = -259887.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1095.0526
$ws.Range("I55").Value = 231.22223
$ws.Range("K55").Value = 231.22223
$ws.Range("M55").Value = -58.22223
$ws.Range("H132").Value = 3267.8823
$ws.Range("I132").Value = 2352.4167
$ws.Range("J132").Value = 5465
$ws.Range("K132").Value = 7057.250100000001
$ws.Range("L132").Value = 16395
$ws.Range("M132").Value = -4527.250100000001
$ws.Range("N132").Value = -21455

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 29416624
$ws.Range("I122").Value = 40004852
$ws.Range("J122").Value = 4878.222
$ws.Range("K122").Value = 120014556
$ws.Range("L122").Value = 14634.666
$ws.Range("M122").Value = -120012106
$ws.Range("N122").Value = -19534.666
$ws.Range("H132").Value = 18504.191
$ws.Range("I132").Value = 2258.262
$ws.Range("K132").Value = 6774.786
$ws.Range("M132").Value = -4244.786
$ws.Range("H136").Value = 247667.92
$ws.Range("I136").Value = 273017.3
$ws.Range("K136").Value = 819051.8999999999
$ws.Range("M136").Value = -816501.8999999999
